# Commit: "adding new progress as of date 04 nov 2025"
#
# The "Training Dashboard" sheet tracks, per training/SOP row, how many
# days remain before expiry ("PERIOD TO EXPIRE", column H) as measured
# from the "LAST UPDATE" date (column I). Moving the progress/update date
# forward by one day (03-Nov-2025 -> 04-Nov-2025) shortens every
# remaining-days count in column H by exactly 1, for every data row
# (rows 3 through 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$firstRow = 3
$lastRow  = 14
$periodCol = 8   # H - PERIOD TO EXPIRE
$updateCol = 9   # I - LAST UPDATE

$newUpdateDate = "04-Nov-2025"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    # Column H holds a plain number of days; read it back with Value2
    # (avoids this host's quirky Value getter) and write one less.
    $periodCell = $ws.Cells.Item($row, $periodCol)
    $periodCell.Value = $periodCell.Value2 - 1

    # Column I stores the date as literal text ("DD-MMM-YYYY"), not a
    # real date serial. A leading apostrophe forces the new value to be
    # kept as text instead of being auto-converted to a date value.
    $updateCell = $ws.Cells.Item($row, $updateCol)
    $updateCell.Value = "'" + $newUpdateDate
}
